# Mise à jour de certains champs de Modules et de Professeurs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: swap the last two column headers.
#   C1: "Enseignant"       -> "Chef  Module"  (note: two spaces, matches source)
#   D1: "Nombre d'heures"  -> "Composants"
$ws.Range("C1").Value = "Chef  Module"
$ws.Range("D1").Value = "Composants"

# Column widths for the now-wider header text.
#   Column C (3) -> character width 35
#   Column D (4) -> character width ~24.57 (closest reproducible value: 24.5)
$ws.Columns.Item(3).ColumnWidth = 34.17
$ws.Columns.Item(4).ColumnWidth = 23.67

# Move the active selection from D2 to E8.
$ws.Range("E8").Select() | Out-Null
